# Update the cached "auto date" placeholder text (2022/9/18 -> 2022/9/19)
# across every slide, the slide master, every slide layout, and the notes
# master; also refresh the YouTube link/text on slide 6 to the new video.

$p = $ppt.ActivePresentation

$OLD_DATE = "2022/9/18"
$NEW_DATE = "2022/9/19"

function Update-DateShapes($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)
        if ($shape.HasTextFrame -ne 0) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text -eq $OLD_DATE) {
                $tr.Text = $NEW_DATE
            }
        }
    }
}

# 1) Every slide's date placeholder.
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    Update-DateShapes($s.Shapes)
}

# 2) The slide master's date placeholder.
$master = $p.SlideMaster
Update-DateShapes($master.Shapes)

# 3) Every slide layout's date placeholder.
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DateShapes($layout.Shapes)
}

# 4) The notes master's date placeholder.
$notesMaster = $p.NotesMaster
Update-DateShapes($notesMaster.Shapes)

# 5) Slide 6: refresh the YouTube link text to the new video/index.
$slide6 = $p.Slides.Item(6)
for ($i = 1; $i -le $slide6.Shapes.Count; $i++) {
    $shape = $slide6.Shapes.Item($i)
    if ($shape.HasTextFrame -ne 0) {
        $tr = $shape.TextFrame.TextRange
        if ($tr.Text -like "*youtube.com/watch?v=lDmYXYEdp7I*") {
            $tr.Text = "https://www.youtube.com/watch?v=a-wc5ldxM7k&list=PL1qVKHVG3ZfVb91esBQ0-0SQC3dGGeXkn&index=4"
        }
    }
}
